$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Test"
$ws.Range("B6").Value = "Write"
$ws.Range("A7").Value = "Test"
$ws.Range("B7").Value = "Write"
$ws.Range("A8").Value = "Test"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "10000"
$ws.Range("B8").Style = "Normal"
$ws.Range("A9").Value = "Test"
$ws.Range("B9").Value = "Demo"

[void]$ws.Range("I12").Select()
